$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the subject / date identifier cells (kept as Text, matching existing column format)
$ws.Range("A2").Value = "Xe-027"
$ws.Range("B2").Value = "2021-11-05"
$ws.Range("C2").Value = "2021-11-15"

# Update the numeric ventilation-summary measurements for the (re-run / right-sized) subject
$ws.Range("D2").Value = 23.18216044855421
$ws.Range("E2").Value = 6.3431807745293129
$ws.Range("F2").Value = 19.014235420174501
$ws.Range("G2").Value = 6.0921475585489056
$ws.Range("H2").Value = 25.357416194703813
$ws.Range("I2").Value = 6.0370427062605234
$ws.Range("J2").Value = 17.624368590234198
$ws.Range("K2").Value = 2.5225776825348234
$ws.Range("L2").Value = 23.66141129649472
$ws.Range("M2").Value = 20.719424460431654
$ws.Range("N2").Value = 32.769018827491195
$ws.Range("O2").Value = 26.738098882596052
$ws.Range("P2").Value = 12.37410071942446
$ws.Range("Q2").Value = 4.8951477116179403
$ws.Range("R2").Value = 2.5042093984386957
$ws.Range("S2").Value = 14.780345936017145
$ws.Range("T2").Value = 17.449869891320986
$ws.Range("U2").Value = 22.045002296035513
$ws.Range("V2").Value = 18.423388948415738
$ws.Range("W2").Value = 11.930200520434717
$ws.Range("X2").Value = 15.371192407775908
$ws.Range("Y2").Value = 5.6696770243379762
$ws.Range("Z2").Value = 28.675952854737485
$ws.Range("AA2").Value = 52.337364151232208
$ws.Range("AB2").Value = 13.31700596969233
